$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240506768226624
$ws.Range("B1").Value = 2.363690853118896
$ws.Range("C1").Value = 3.871894598007202
$ws.Range("D1").Value = 3.141841173171997
$ws.Range("E1").Value = 1.270837426185608
